$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Text Cards to implement" section (column B header, column C list) ---
$ws.Range("B9").Value = "Text Cards to implement:"
$ws.Range("B9").Font.Bold = $true

$ws.Range("C9").Value = "Adventurer"
$ws.Range("C10").Value = "Bureaucrat"
$ws.Range("C11").Value = "Cellar"
$ws.Range("C12").Value = "Chancellor"
$ws.Range("C13").Value = "Chapel"
$ws.Range("C14").Value = "Council Room"
$ws.Range("C15").Value = "Feast"
$ws.Range("C16").Value = "Gardens"
$ws.Range("C17").Value = "Library"
$ws.Range("C18").Value = "Militia"
$ws.Range("C19").Value = "Mine"
$ws.Range("C20").Value = "Moat"
$ws.Range("C21").Value = "Money Lender"
$ws.Range("C22").Value = "Remodel"
$ws.Range("C23").Value = "Spy"
$ws.Range("C24").Value = "Thief"
$ws.Range("C25").Value = "Throne Room"
$ws.Range("C26").Value = "Witch"
$ws.Range("C27").Value = "Workshop"

# --- "Finished Cards" section (column D header, column C list continues) ---
$ws.Range("D9").Value = "Finished Cards:"
$ws.Range("D9").Font.Bold = $true

$ws.Range("C31").Value = "Copper"
$ws.Range("C32").Value = "Silver"
$ws.Range("C33").Value = "Gold"
$ws.Range("C28").Value = "Estate"
$ws.Range("C29").Value = "Duchy"
$ws.Range("C30").Value = "Province"

# --- Additional finished cards laid out in columns D and E ---
$ws.Range("D12").Value = "Festival"
$ws.Range("D13").Value = "Laboratory"
$ws.Range("D14").Value = "Market"

$ws.Range("E12").Value = "Smithy"
$ws.Range("E13").Value = "Village"
$ws.Range("E14").Value = "Woodcutter"

# --- Metrics (percent finished / number finished) ---
$ws.Range("E10").Value = "percent finished"
$ws.Range("E10").Font.Bold = $true

$ws.Range("E9").Value = "number finished"
$ws.Range("E9").Font.Bold = $true

$ws.Range("F9").Value = 6
$ws.Range("F10").Formula = "=(F9 / 31) * 100"

# --- Column E width adjustment ---
$ws.Columns("E").ColumnWidth = 13.43

# --- Selection matches the end state recorded in the workbook ---
[void]$ws.Range("F13").Select()
